# Add two new rows of data to Sheet1, mirroring the existing A/B/C pattern
# (range label in column A, numeric id already in column B, owner name in column C),
# then leave the selection on L6 as in the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "500-599"
$ws.Range("C6").Value = "Andrew"

$ws.Range("A7").Value = "600-699"
$ws.Range("C7").Value = "Miho"

$ws.Range("L6").Select()
